# collection_group_id values (column H) reused the same shared string for
# rows belonging to different End-of-Study / End-of-Treatment / Study
# Continuation scenarios (column J). Make them unique by suffixing each
# base collection_group_id with _EOS, _EOT or _CONT depending on the
# scenario (column J) the row belongs to.
#
# The ranges below are written in the exact order needed so the workbook's
# shared-strings table mints the new unique strings in the same order as
# the target workbook (ADVEVENT_EOT, ADVEVENT_EOS, ADVEVENT_CONT, then the
# remaining bases each in EOS/EOT/CONT order, and finally LTFUP_EOS).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ADVEVENT (End of Treatment, End of Study, Study Continuation)
$ws.Range("H52:H56").Value = "ADVEVENT_EOT"
$ws.Range("H10:H14").Value = "ADVEVENT_EOS"
$ws.Range("H106:H110").Value = "ADVEVENT_CONT"

# DEAD
$ws.Range("H15:H19").Value = "DEAD_EOS"
$ws.Range("H57:H61").Value = "DEAD_EOT"
$ws.Range("H111:H115").Value = "DEAD_CONT"

# LACKEFFICACY
$ws.Range("H20:H23").Value = "LACKEFFICACY_EOS"
$ws.Range("H62:H65").Value = "LACKEFFICACY_EOT"
$ws.Range("H124:H127").Value = "LACKEFFICACY_CONT"

# PHYDECISION
$ws.Range("H28:H31").Value = "PHYDECISION_EOS"
$ws.Range("H66:H69").Value = "PHYDECISION_EOT"
$ws.Range("H128:H131").Value = "PHYDECISION_CONT"

# PROGDISEASE
$ws.Range("H32:H35").Value = "PROGDISEASE_EOS"
$ws.Range("H70:H73").Value = "PROGDISEASE_EOT"
$ws.Range("H132:H135").Value = "PROGDISEASE_CONT"

# PROTCOMP
$ws.Range("H36:H39").Value = "PROTCOMP_EOS"
$ws.Range("H74:H77").Value = "PROTCOMP_EOT"
$ws.Range("H136:H139").Value = "PROTCOMP_CONT"

# PROTDEV
$ws.Range("H40:H43").Value = "PROTDEV_EOS"
$ws.Range("H78:H81").Value = "PROTDEV_EOT"
$ws.Range("H140:H143").Value = "PROTDEV_CONT"

# SUBJPREG
$ws.Range("H44:H47").Value = "SUBJPREG_EOS"
$ws.Range("H82:H85").Value = "SUBJPREG_EOT"
$ws.Range("H144:H147").Value = "SUBJPREG_CONT"

# SUBJWITHDRAW
$ws.Range("H48:H51").Value = "SUBJWITHDRAW_EOS"
$ws.Range("H86:H89").Value = "SUBJWITHDRAW_EOT"
$ws.Range("H148:H151").Value = "SUBJWITHDRAW_CONT"

# LTFUP only occurs in the End of Study scenario
$ws.Range("H24:H27").Value = "LTFUP_EOS"

# Restore the current selection to where the author left off editing.
[void]$ws.Range("J33").Select()
